# "Season up to 1/17" update:
#  - The game previously listed as the next upcoming game on the "Next"
#    sheet (2024-01-15 vs NOP, home) has now been played. Its result is
#    appended as a new row to the "Games" sheet, and the corresponding
#    row is removed from the "Next" sheet (all following rows shift up).

$wb = $excel.ActiveWorkbook

$games = $wb.Worksheets.Item("Games")
$next  = $wb.Worksheets.Item("Next")

# --- 1) Append the completed game to the "Games" sheet (new row 42) ---
$newRow = $games.UsedRange.Rows.Count + 1

$games.Cells.Item($newRow, 1).Value  = 41            # Game
$games.Cells.Item($newRow, 2).Value  = 45306          # Date
$games.Cells.Item($newRow, 2).NumberFormat = "YYYY-MM-DD"
$games.Cells.Item($newRow, 3).Value  = 1              # Streak
$games.Cells.Item($newRow, 4).Value  = 125            # Pts
$games.Cells.Item($newRow, 5).Value  = 96             # Pace
$games.Cells.Item($newRow, 6).Value  = 0.565           # eFG
$games.Cells.Item($newRow, 7).Value  = 9.9            # TOV
$games.Cells.Item($newRow, 8).Value  = 31             # ORB
$games.Cells.Item($newRow, 9).Value  = 0.341          # FTR
$games.Cells.Item($newRow, 10).Value = 130.3          # ORT
$games.Cells.Item($newRow, 11).Value = "NOP"          # OppID
$games.Cells.Item($newRow, 12).Value = 120            # OppPts
$games.Cells.Item($newRow, 13).Value = 0.549          # OppeFG
$games.Cells.Item($newRow, 14).Value = 8.3            # OppTOV
$games.Cells.Item($newRow, 15).Value = 21.6           # OppORB
$games.Cells.Item($newRow, 16).Value = 0.383          # OppFTR
$games.Cells.Item($newRow, 17).Value = 125.1          # OppORT
$games.Cells.Item($newRow, 18).Value = 1              # Location
$games.Cells.Item($newRow, 19).Value = 1              # Target

# --- 2) Remove that game from the "Next" sheet; remaining rows shift up ---
$next.Rows.Item(2).Delete()
